$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.388317276604681
$ws.Range("C2").Value = 0.04981286588794376
$ws.Range("D2").Value = 0.1577979911730782
$ws.Range("E2").Value = 0.07019828564042463
$ws.Range("F2").Value = 2.625751773257235
$ws.Range("K2").Value = 0.9305988667144049
$ws.Range("L2").Value = 0.1963251175074276
$ws.Range("M2").Value = 0.2988679867042947
$ws.Range("N2").Value = 3.740235092271703
$ws.Range("B3").Value = 1.347845959518111
$ws.Range("C3").Value = 0.04514965112367975
$ws.Range("D3").Value = 0.1581448829327385
$ws.Range("E3").Value = 0.07033141536890586
$ws.Range("F3").Value = 2.602581907047437
$ws.Range("K3").Value = 0.8883912467549067
$ws.Range("L3").Value = 0.1940149243992479
$ws.Range("M3").Value = 0.2919269069122734
$ws.Range("N3").Value = 3.742305222211286
$ws.Range("B4").Value = 1.323817732215304
$ws.Range("C4").Value = 0.04226370416463965
$ws.Range("D4").Value = 0.1583614021332123
$ws.Range("E4").Value = 0.07043031795492904
$ws.Range("F4").Value = 2.589425679108871
$ws.Range("K4").Value = 0.8630588559734633
$ws.Range("L4").Value = 0.1926930427625848
$ws.Range("M4").Value = 0.2878346556716167
$ws.Range("N4").Value = 3.74434281467633
$ws.Range("B5").Value = 1.314232525789095
$ws.Range("C5").Value = 0.04108186791624746
$ws.Range("D5").Value = 0.1584505233951292
$ws.Range("E5").Value = 0.07047494563249934
$ws.Range("F5").Value = 2.584333303512096
$ws.Range("K5").Value = 0.8528821651779026
$ws.Range("L5").Value = 0.1921786667684984
$ws.Range("M5").Value = 0.2862096803863245
$ws.Range("N5").Value = 3.745365747169942
$ws.Range("B6").Value = 1.312653380257643
$ws.Range("C6").Value = 0.04088527297363953
$ws.Range("D6").Value = 0.1584653755781078
$ws.Range("E6").Value = 0.07048261745186224
$ws.Range("F6").Value = 2.583503953079244
$ws.Range("K6").Value = 0.8512011738118304
$ws.Range("L6").Value = 0.1920947232103103
$ws.Range("M6").Value = 0.2859424305853508
$ws.Range("N6").Value = 3.745547232352408
$ws.Range("B7").Value = 1.323687626740679
$ws.Range("C7").Value = 0.04224778905044957
$ws.Range("D7").Value = 0.1583626004580125
$ws.Range("E7").Value = 0.07043090229926818
$ws.Range("F7").Value = 2.589355913093144
$ws.Range("K7").Value = 0.862921016796264
$ws.Range("L7").Value = 0.1926860072924157
$ws.Range("M7").Value = 0.2878125679921517
$ws.Range("N7").Value = 3.744355830705857
$ws.Range("B8").Value = 1.374192259284939
$ws.Range("C8").Value = 0.04820965574612046
$ws.Range("D8").Value = 0.1579168683521264
$ws.Range("E8").Value = 0.07024063242698642
$ws.Range("F8").Value = 2.617540552107769
$ws.Range("K8").Value = 0.9159244858464319
$ws.Range("L8").Value = 0.1955085235183276
$ws.Range("M8").Value = 0.2964395004089369
$ws.Range("N8").Value = 3.740789685325154
$ws.Range("B9").Value = 1.479757675600297
$ws.Range("C9").Value = 0.05972461984342203
$ws.Range("D9").Value = 0.1570706872285257
$ws.Range("E9").Value = 0.07000330739512428
$ws.Range("F9").Value = 2.681316881054997
$ws.Range("K9").Value = 1.024507976926998
$ws.Range("L9").Value = 0.201809773882232
$ws.Range("M9").Value = 0.314703810579104
$ws.Range("N9").Value = 3.739887086964032
$ws.Range("B10").Value = 1.561318408347347
$ws.Range("C10").Value = 0.0680834397310548
$ws.Range("D10").Value = 0.1564658558215246
$ws.Range("E10").Value = 0.06991127517283147
$ws.Range("F10").Value = 2.733388172271404
$ws.Range("K10").Value = 1.107147313520244
$ws.Range("L10").Value = 0.2069072007139994
$ws.Range("M10").Value = 0.3289473801535507
$ws.Range("N10").Value = 3.742951321020342
$ws.Range("B11").Value = 1.599297271945261
$ws.Range("C11").Value = 0.07186548398264847
$ws.Range("D11").Value = 0.1561943312605623
$ws.Range("E11").Value = 0.0698871926379816
$ws.Range("F11").Value = 2.758215998717247
$ws.Range("K11").Value = 1.145371365710162
$ws.Range("L11").Value = 0.20932797363524
$ws.Range("M11").Value = 0.3356071743520559
$ws.Range("N11").Value = 3.745158021733431
$ws.Range("B12").Value = 1.613805177091479
$ws.Range("C12").Value = 0.07329480452816028
$ws.Range("D12").Value = 0.1560920308511262
$ws.Range("E12").Value = 0.06988062239071802
$ws.Range("F12").Value = 2.767782036644519
$ws.Range("K12").Value = 1.159936972332872
$ws.Range("L12").Value = 0.2102593172547387
$ws.Range("M12").Value = 0.3381550330237673
$ws.Range("N12").Value = 3.746110755368392
$ws.Range("B13").Value = 1.610675028672347
$ws.Range("C13").Value = 0.07298710059392022
$ws.Range("D13").Value = 0.1561140399426968
$ws.Range("E13").Value = 0.06988192415659888
$ws.Range("F13").Value = 2.765714508011769
$ws.Range("K13").Value = 1.15679595758084
$ws.Range("L13").Value = 0.2100580843188879
$ws.Range("M13").Value = 0.3376051522638193
$ws.Range("N13").Value = 3.745900355196966
$ws.Range("B14").Value = 1.600488317540339
$ws.Range("C14").Value = 0.07198313185689642
$ws.Range("D14").Value = 0.1561859045217835
$ws.Range("E14").Value = 0.06988660105349886
$ws.Range("F14").Value = 2.758999708277571
$ws.Range("K14").Value = 1.146567861948938
$ws.Range("L14").Value = 0.2094043022798076
$ws.Range("M14").Value = 0.3358162683029846
$ws.Range("N14").Value = 3.745234055436669
$ws.Range("B15").Value = 1.59426509219611
$ws.Range("C15").Value = 0.0713678022483748
$ws.Range("D15").Value = 0.156229991380691
$ws.Range("E15").Value = 0.06988979754584967
$ws.Range("F15").Value = 2.754908102541805
$ws.Range("K15").Value = 1.140314714856714
$ws.Range("L15").Value = 0.2090057494784787
$ws.Range("M15").Value = 0.3347239035558474
$ws.Range("N15").Value = 3.744841184867596
$ws.Range("B16").Value = 1.55885404540652
$ws.Range("C16").Value = 0.06783587003471325
$ws.Range("D16").Value = 0.1564836733577142
$ws.Range("E16").Value = 0.06991320618937458
$ws.Range("F16").Value = 2.731788587368868
$ws.Range("K16").Value = 1.104662000445302
$ws.Range("L16").Value = 0.2067510478909327
$ws.Range("M16").Value = 0.3285157762293451
$ws.Range("N16").Value = 3.742823482351483
$ws.Range("B17").Value = 1.537355068878242
$ws.Range("C17").Value = 0.0656639691968337
$ws.Range("D17").Value = 0.1566402252233132
$ws.Range("E17").Value = 0.06993211570713953
$ws.Range("F17").Value = 2.717897792047296
$ws.Range("K17").Value = 1.082952005912006
$ws.Range("L17").Value = 0.2053939620204233
$ws.Range("M17").Value = 0.3247534821263187
$ws.Range("N17").Value = 3.741794017593776
$ws.Range("B18").Value = 1.525071918107585
$ws.Range("C18").Value = 0.06441282201663512
$ws.Range("D18").Value = 0.1567306100855834
$ws.Range("E18").Value = 0.06994466628778007
$ws.Range("F18").Value = 2.710015489341146
$ws.Range("K18").Value = 1.07052436894395
$ws.Range("L18").Value = 0.2046229964754644
$ws.Range("M18").Value = 0.3226064822985606
$ws.Range("N18").Value = 3.741278384042687
$ws.Range("B19").Value = 1.520927219525504
$ws.Range("C19").Value = 0.06398887165376266
$ws.Range("D19").Value = 0.1567612713076283
$ws.Range("E19").Value = 0.06994920350917333
$ws.Range("F19").Value = 2.70736509955735
$ws.Range("K19").Value = 1.066326775262411
$ws.Range("L19").Value = 0.2043636088723133
$ws.Range("M19").Value = 0.3218824597948924
$ws.Range("N19").Value = 3.741116929230259
$ws.Range("B20").Value = 1.539635133361628
$ws.Range("C20").Value = 0.06589537063551631
$ws.Range("D20").Value = 0.1566235247845302
$ws.Range("E20").Value = 0.06992992952792676
$ws.Range("F20").Value = 2.719365382071416
$ws.Range("K20").Value = 1.085256924740634
$ws.Range("L20").Value = 0.2055374332066862
$ws.Range("M20").Value = 0.3251522281362966
$ws.Range("N20").Value = 3.741895688121872
$ws.Range("B21").Value = 1.603476977481137
$ws.Range("C21").Value = 0.07227809872782132
$ws.Range("D21").Value = 0.1561647820377381
$ws.Range("E21").Value = 0.06988515821152497
$ws.Range("F21").Value = 2.760967546818762
$ws.Range("K21").Value = 1.149569630164279
$ws.Range("L21").Value = 0.2095959363887516
$ws.Range("M21").Value = 0.3363410028262805
$ws.Range("N21").Value = 3.745426583808296
$ws.Range("B22").Value = 1.645936567904585
$ws.Range("C22").Value = 0.07643296182924075
$ws.Range("D22").Value = 0.1558679985000175
$ws.Range("E22").Value = 0.06987075313214852
$ws.Range("F22").Value = 2.78911465389325
$ws.Range("K22").Value = 1.192132167361478
$ws.Range("L22").Value = 0.2123337879645248
$ws.Range("M22").Value = 0.3438047122479944
$ws.Range("N22").Value = 3.7484169094609
$ws.Range("B23").Value = 1.623207775367632
$ws.Range("C23").Value = 0.07421692742076402
$ws.Range("D23").Value = 0.1560261200400852
$ws.Range("E23").Value = 0.06987708471121934
$ws.Range("F23").Value = 2.774004283645354
$ws.Range("K23").Value = 1.169367118129486
$ws.Range("L23").Value = 0.2108647346633461
$ws.Range("M23").Value = 0.3398073520782887
$ws.Range("N23").Value = 3.746758374091371
$ws.Range("B24").Value = 1.538604076693446
$ws.Range("C24").Value = 0.06579076181574806
$ws.Range("D24").Value = 0.1566310738649452
$ws.Range("E24").Value = 0.0699309126685268
$ws.Range("F24").Value = 2.718701561623433
$ws.Range("K24").Value = 1.084214703604516
$ws.Range("L24").Value = 0.2054725410986293
$ws.Range("M24").Value = 0.3249719052498961
$ws.Range("N24").Value = 3.741849485464172
$ws.Range("B25").Value = 1.450498390351697
$ws.Range("C25").Value = 0.05662774336221332
$ws.Range("D25").Value = 0.1572966310018504
$ws.Range("E25").Value = 0.07005301959849142
$ws.Range("F25").Value = 2.663150046952779
$ws.Range("K25").Value = 0.994632822831278
$ws.Range("L25").Value = 0.2000230065989754
$ws.Range("M25").Value = 0.3096182602905984
$ws.Range("N25").Value = 3.739477753214175
